# Apply the "flash分配" workbook fixes described in the commit message:
#   1. onboard led en/disable                      (no visible change needed here in sheet data)
#   2. load parameter flash save and restore        -> add a new "load_radar_parameter" row,
#                                                        and clear the stale "factory value" column
#
# Concretely, on Sheet1:
#   - Clear the (now-invalid / to-be-recomputed) factory values in G2:G9, keeping their formatting.
#   - Insert a brand-new row 16 documenting the new `load_radar_parameter` flash variable
#     (u32 @ UPSSA0 offset 0x38), pushing every row from the old 16 down by one.
#   - Update the sheet's cursor/selection like a user would leave it after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Clear the stale factory-default numbers in column G (rows 2-9) ---------------
# ClearContents keeps the existing cell style/border/number-format, only removes the value.
$ws.Range("G2:G9").ClearContents() | Out-Null

# --- 2. Insert the new "load_radar_parameter" row above the current row 16 -----------
$ws.Rows.Item(16).Insert() | Out-Null

# Pick up the formatting of the (now shifted-down) row that used to be row 16 so the
# freshly inserted row looks consistent with its neighbours (border/alignment/number format).
$ws.Range("A17:G17").Copy() | Out-Null
$ws.Range("A16:G16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new row's content.
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "load_radar_parameter"
$ws.Range("C16").Value = "u32"
$ws.Range("D16").Value = "加载覆盖范围"
$ws.Range("E16").Value = "UPSSA0"
$ws.Range("F16").Value = "0x38"
$ws.Range("G16").Value = 1

# --- 3. Leave the selection where the editing user ended up -------------------------
$ws.Range("C24").Select() | Out-Null
